$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_G_acc_LT"
$ws.Range("A2").Value = 86.141304347826093
$ws.Range("A3").Value = 85.326086956521735
$ws.Range("A4").Value = 86.41304347826086
$ws.Range("A5").Value = 78.804347826086953
$ws.Range("A6").Value = 79.891304347826093
$ws.Range("A7").Value = 81.25
$ws.Range("A8").Value = 89.945652173913047
$ws.Range("A9").Value = 88.858695652173907
$ws.Range("A10").Value = 88.58695652173914
$ws.Range("A11").Value = 88.858695652173907
$ws.Range("A12").Value = 73.91304347826086
$ws.Range("A13").Value = 80.706521739130437
$ws.Range("A14").Value = 87.228260869565219
$ws.Range("A15").Value = 86.41304347826086
$ws.Range("A16").Value = 86.956521739130437
$ws.Range("A17").Value = 74.728260869565219
$ws.Range("A18").Value = 76.902173913043484
$ws.Range("A19").Value = 79.619565217391312
$ws.Range("A20").Value = 87.228260869565219
$ws.Range("A21").Value = 88.58695652173914
$ws.Range("A22").Value = 88.315217391304344
$ws.Range("A23").Value = 76.08695652173914
$ws.Range("A24").Value = 77.173913043478265
$ws.Range("A25").Value = 75.815217391304344
$ws.Range("A26").Value = 86.141304347826093
$ws.Range("A27").Value = 83.967391304347828
$ws.Range("A28").Value = 83.152173913043484
$ws.Range("A29").Value = 82.608695652173907
$ws.Range("A30").Value = 79.891304347826093
$ws.Range("A31").Value = 83.152173913043484
$ws.Range("A32").Value = 81.521739130434781
$ws.Range("A33").Value = 83.423913043478265
$ws.Range("A34").Value = 82.065217391304344
$ws.Range("A35").Value = 80.706521739130437
$ws.Range("A36").Value = 81.25
$ws.Range("A37").Value = 73.91304347826086
$ws.Range("A38").Value = 81.521739130434781
$ws.Range("A39").Value = 79.076086956521735
$ws.Range("A40").Value = 79.076086956521735
$ws.Range("A41").Value = 83.967391304347828
$ws.Range("A42").Value = 86.41304347826086
$ws.Range("A43").Value = 84.782608695652172
$ws.Range("A44").Value = 85.054347826086953
$ws.Range("A45").Value = 86.141304347826093
$ws.Range("A46").Value = 85.869565217391312
$ws.Range("A47").Value = 79.076086956521735
$ws.Range("A48").Value = 75
$ws.Range("A49").Value = 82.608695652173907
